$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update designator references following the capacitor renumbering ---
# C1,C2 shared part -> now only C1 uses it
$ws.Range("B5").Value = "C1"
# C3 alone -> now C2,C3 share this part
$ws.Range("B8").Value = "C2,C3"

# --- Row 9: old R4 (470R resistor) replaced with new U2 connector part ---
$ws.Range("A9").Value = "0.4mm 2 24P Brick nogging Female SMD,P=0.4mm Mezzanine Connectors (Board to Board) ROHS"
$ws.Range("B9").Value = "U2"
$ws.Range("C9").Value = "SMD,P=0.4mm"
$ws.Range("D9").Value = "C3640874"
# D9 previously used the wrap-text variant (needed for the old long resistor
# comment); the new short part number doesn't need wrapping, so align its
# formatting with the rest of the row (matches B9/C9's plain style).
$ws.Range("C9").Copy()
$ws.Range("D9").PasteSpecial(-4122)
# Row no longer needs the taller wrapped-text height once D9 stops using the wrap style
$ws.Rows.Item(9).RowHeight = 13.8

# --- Row 10: old R5 (680R resistor) removed, row cleared but kept (4 blank cells) ---
$ws.Range("A10:D10").ClearContents()

# --- Row 11: old D1 (blue LED) removed entirely except col A remains as an empty placeholder ---
$ws.Range("B11:D11").Clear()
$ws.Range("A11").ClearContents()

# --- Row 12: old D2 (yellow LED) removed, row cleared but kept (4 blank cells) ---
$ws.Range("A12:D12").ClearContents()

# --- Row 13: already blank row, reduced down to just col A ---
$ws.Range("B13:D13").Clear()

# --- Row 14: stays as a blank 4-cell row ---
$ws.Range("A14:D14").ClearContents()

# --- Update the active selection, matching the end of the user's edit session ---
$ws.Range("A14").Select()
